# Generate Report for Handback
# Adds "Latest Target File" (F) and "Latest Handback File" (G) hyperlinked
# values for both data rows on the zh-cn and de-de sheets, updates the
# Status column to reflect a completed handback, and sets the
# "Latest Handback DateTime" values.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # RGB(0x64,0x95,0xED) == the workbook's HyperLink font color
$hyperlinkUnderline = 2     # xlUnderlineStyleSingle

function Style-AsHyperlink($range) {
    $range.Font.Underline = $hyperlinkUnderline
    $range.Font.Color = $hyperlinkColor
}

$sheetsInfo = @(
    @{
        Name = "zh-cn"
        MdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/a6166e952e59def1915e5c6d4077ec99dc5cb9d6/e2e/1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848.md"
        MdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/a6166e952e59def1915e5c6d4077ec99dc5cb9d6/e2e/ffff1f6804c2-9b89-4dcb-86c9-d98c33ec6e5e.md"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/300c63d3d18469202d6fcfe76bf5fc9df2c6128c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848.08af7258856ed5122eb31db14c922250c28858c4.zh-cn.xlf"
        XlfName = "1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848.08af7258856ed5122eb31db14c922250c28858c4.zh-cn.xlf"
        HandbackDateTime = "2016-03-18 08:47:32"
    },
    @{
        Name = "de-de"
        MdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/a6166e952e59def1915e5c6d4077ec99dc5cb9d6/e2e/1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848.md"
        MdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/a6166e952e59def1915e5c6d4077ec99dc5cb9d6/e2e/ffff1f6804c2-9b89-4dcb-86c9-d98c33ec6e5e.md"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e25c93479bb96fcc0f4e0728cc736f444a250a2a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848.08af7258856ed5122eb31db14c922250c28858c4.de-de.xlf"
        XlfName = "1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848.08af7258856ed5122eb31db14c922250c28858c4.de-de.xlf"
        HandbackDateTime = "2016-03-18 08:47:37"
    }
)

$mdDisplay = "1f641bcb-6eb4-4a1c-8351-6d8c6ffc8848.md"
$newStatus = "Handed back: in sync with en-US"

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Row 2 - Latest Target File (F) / Latest Handback File (G)
    $ws.Hyperlinks.Add($ws.Range("F2"), $info.MdUrl1, "", "", $mdDisplay) | Out-Null
    Style-AsHyperlink $ws.Range("F2")

    $ws.Hyperlinks.Add($ws.Range("G2"), $info.XlfUrl, "", "", $info.XlfName) | Out-Null
    Style-AsHyperlink $ws.Range("G2")

    # Row 3 - Latest Target File (F) / Latest Handback File (G)
    $ws.Hyperlinks.Add($ws.Range("F3"), $info.MdUrl2, "", "", $mdDisplay) | Out-Null
    Style-AsHyperlink $ws.Range("F3")

    $ws.Hyperlinks.Add($ws.Range("G3"), $info.XlfUrl, "", "", $info.XlfName) | Out-Null
    Style-AsHyperlink $ws.Range("G3")

    # Status column (C) -> handed back
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Handback DateTime (H)
    $ws.Range("H2").Value = $info.HandbackDateTime
    $ws.Range("H3").Value = $info.HandbackDateTime
}

Write-Host "Handback report generated."
